$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.847.76"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.134.70"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.08%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.127.37"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "3.649.44"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "63.836.95"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "3.132.80"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.74%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  -6.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -7.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  -6.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0393"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  -10.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "427.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.62%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "2.871.19"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("E47").Value = "  -7.10%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
